$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: un-merge D2:O2, drop the centered alignment (keep vertical-center + wrap), add Q2 ---
$ws.Range("D2:O2").UnMerge()
$ws.Range("D2:O2").HorizontalAlignment = 1  # xlHAlignGeneral -> removes "horizontal=center", keeps vertical/wrap

$ws.Range("P2").Copy()
$ws.Range("Q2").PasteSpecial(-4122)  # xlPasteFormats (same border/font as P2)
$ws.Range("Q2").Interior.Pattern = -4142  # xlPatternNone (Q2 carries no fill, unlike P2)

$ws.Rows("2:2").RowHeight = 15

# --- Row 1: A1 style index shifts (no visual change) because an old unused xf got dropped ---
# (handled implicitly by the cellXfs table produced by the other edits below)

# --- Row 3: new header year 2020 in Q3, same look as P3 ---
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q3").Value = 2020

# --- Row 4: new data point in Q4, same look as P4 ---
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q4").Value = 14.5

# --- Row 5: new data point in Q5, same look as P5 ---
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q5").Value = 13.8

$ws.Application.CutCopyMode = $false

# --- sheet view: selection moves to P13 ---
$ws.Range("P13").Select()

Write-Host "done"
